$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date serial that was bumped by one day
# (2026-02-07 -> 2026-02-08, i.e. Excel serial 46060 -> 46061) for every data
# row (rows 2 through 329).
$ws.Range("C2:C329").Value = 46061
